$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D/E columns for rows with new price/volume values ---
$ws.Range("D2").Value = "59.372.15"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "2.605.70"
$ws.Range("E3").Value = "  +1.77%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "536.73"
$ws.Range("E5").Value = "  +4.53%  "

$ws.Range("D6").Value = "140.90"
$ws.Range("E6").Value = "  +3.30%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "0.569"
$ws.Range("E8").Value = "  +1.78%  "

$ws.Range("D9").Value = "2.618.17"
$ws.Range("E9").Value = "  +1.74%  "

$ws.Range("D10").Value = "6.47"
$ws.Range("E10").Value = "  -0.42%  "

$ws.Range("E11").Value = "  +4.61%  "

$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +3.66%  "

$ws.Range("E13").Value = "  +2.09%  "

$ws.Range("D14").Value = "3.068.70"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").Value = "59.307.62"
$ws.Range("E15").Value = "  +1.85%  "

$ws.Range("D16").Value = "20.56"
$ws.Range("E16").Value = "  +2.42%  "

$ws.Range("D17").Value = "2.607.44"
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").Value = "345.92"
$ws.Range("E19").Value = "  +3.83%  "

$ws.Range("E20").Value = "  +2.27%  "

$ws.Range("D21").Value = "10.14"
$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "67.09"
$ws.Range("E24").Value = "  +2.20%  "

$ws.Range("D25").Value = "0.167"
$ws.Range("E25").Value = "  +1.42%  "

$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").Value = "7.21"
$ws.Range("E28").Value = "  +4.27%  "

$ws.Range("D29").Value = "0.0₃0748"
$ws.Range("E29").Value = "  +7.96%  "

$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +6.37%  "

$ws.Range("D32").Value = "5.87"
$ws.Range("E32").Value = "  +1.07%  "

$ws.Range("D33").Value = "18.91"
$ws.Range("E33").Value = "  +2.05%  "

$ws.Range("D34").Value = "149.18"
$ws.Range("E34").Value = "  +0.63%  "

$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +3.62%  "

$ws.Range("D36").Value = "1.12"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").Value = "36.93"
$ws.Range("E37").Value = "  +1.84%  "

$ws.Range("D38").Value = "0.841"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("E39").Value = "  +4.13%  "

$ws.Range("D40").Value = "0.838"
$ws.Range("E40").Value = "  +3.56%  "

$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("D44").Value = "0.600"
$ws.Range("E44").Value = "  +3.02%  "

$ws.Range("D47").Value = "0.0524"
$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").Value = "1.948.09"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").Value = "0.0223"
$ws.Range("E49").Value = "  +3.34%  "

$ws.Range("D50").Value = "18.35"
$ws.Range("E50").Value = "  +5.67%  "

$ws.Range("D51").Value = "4.52"
$ws.Range("E51").Value = "  +3.31%  "

# --- Row 42/43 swap: Bittensor <-> FirstDigitalUSD ---
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "276.44"
$ws.Range("E43").Value = "  +3.49%  "

# --- Row 45/46 swap: Stellar <-> WhiteBITCoin ---
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.75"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.0962"
$ws.Range("E46").Value = "  +2.57%  "
